# Add the new "2022-Q4" detail sheet, inserted right after "总计" and
# right before the existing "2022-Q3" sheet.
$wb = $excel.ActiveWorkbook

$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($beforeSheet)
$q4.Name = "2022-Q4"

# Header row (same layout/style as the other quarterly detail sheets).
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4Data = @(
    @(0, "506005", "博时科创板三年定开混合",             "20.34", "98.81", "4.93", "1.0028", 5),
    @(1, "410003", "华富成长趋势混合",                   "10.74", "87.35", "4.45", "0.4779", 9),
    @(2, "410007", "华富价值增长混合",                   "8.09",  "79.81", "5.28", "0.4272", 6),
    @(3, "014024", "华富卓越成长一年持有期混合A",        "7.94",  "93.53", "4.11", "0.3263", 10),
    @(4, "506007", "广发科创板两年定开混合",             "5.31",  "88.81", "3.66", "0.1943", 10),
    @(5, "009398", "华富成长企业精选股票",               "3.12",  "94.55", "4.45", "0.1388", 10),
    @(6, "014706", "华富匠心明选一年持有期混合A",        "1.93",  "89.78", "4.23", "0.0816", 9),
    @(7, "014707", "华富匠心明选一年持有期混合C",        "1.48",  "89.78", "4.23", "0.0626", 9),
    @(8, "003152", "华富天鑫灵活配置混合A",              "0.93",  "88.63", "4.30", "0.0400", 9),
    @(9, "014025", "华富卓越成长一年持有期混合C",        "0.22",  "93.53", "4.11", "0.0090", 10),
    @(10, "003153", "华富天鑫灵活配置混合C",             "0.16",  "88.63", "4.30", "0.0069", 9)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Match the look of the sibling quarter sheets: bold/bordered header row
# style (column A header cell is blank but still carries the style) plus
# the bold/bordered style used on each data row's index cell (column A).
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2").Copy()
$q4.Range("A2:A12").PasteSpecial(-4122)

# Restore the original "active sheet" (the last tab, "2021-Q3") so the
# workbook's selection state is unchanged by inserting the new tab.
$wb.Worksheets.Item("2021-Q3").Activate()
$q4.Range("A1").Select()

# ---------------------------------------------------------------------------
# Update the "总计" (summary) sheet: insert a new row for 2022-Q4 ahead of
# the 2022-Q3 row, and append the row for 2021-Q3 that now comes back into
# view at the bottom of the table.
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 11
$total.Cells.Item(2, 4).Value = 2.77

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = "2021-Q3"
$total.Cells.Item(7, 3).Value = 4
$total.Cells.Item(7, 4).Value = 0.01
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

$total.Range("A1").Select()
$wb.Worksheets.Item("2021-Q3").Activate()
